$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data in final row order: (month label, value) for A and D columns.
$data = @(
    ,@("2018-10", 96.90000000000001)
    ,@("2018-11", 97.7)
    ,@("2018-12", 98)
    ,@("2018-01", 93.0335)
    ,@("2018-02", 92.7)
    ,@("2018-03", 94.09999999999999)
    ,@("2018-04", 95.09999999999999)
    ,@("2018-05", 93.8)
    ,@("2018-06", 94.40000000000001)
    ,@("2018-07", 95.40000000000001)
    ,@("2018-08", 95.8)
    ,@("2018-09", 96.40000000000001)
    ,@("2019-10", 102.3)
    ,@("2019-11", 102)
    ,@("2019-12", 102.2)
    ,@("2019-01", 102.5)
    ,@("2019-02", 100.6)
    ,@("2019-03", 99.59999999999999)
    ,@("2019-04", 99.40000000000001)
    ,@("2019-05", 101.8)
    ,@("2019-06", 102.4)
    ,@("2019-07", 101.6)
    ,@("2019-08", 101.4)
    ,@("2019-09", 101.4)
    ,@("2020-10", 97.5)
    ,@("2020-11", 96.90000000000001)
    ,@("2020-12", 96.7)
    ,@("2020-01", 99.2)
    ,@("2020-02", 100.4)
    ,@("2020-03", 101.7)
    ,@("2020-04", 100.6)
    ,@("2020-05", 100)
    ,@("2020-06", 98.90000000000001)
    ,@("2020-07", 99.8)
    ,@("2020-08", 99.2)
    ,@("2020-09", 98.8)
    ,@("2021-10", 97.90000000000001)
    ,@("2021-11", 97.8)
    ,@("2021-12", 98.09999999999999)
    ,@("2021-01", 99.40000000000001)
    ,@("2021-02", 98.2)
    ,@("2021-03", 98.2)
    ,@("2021-04", 97.40000000000001)
    ,@("2021-05", 96.7)
    ,@("2021-06", 96.8)
    ,@("2021-07", 101.6)
    ,@("2021-08", 101.6)
    ,@("2021-09", 97.3)
    ,@("2022-10", 100.6)
    ,@("2022-11", 100.9)
    ,@("2022-12", 100.8)
    ,@("2022-01", 98.40000000000001)
    ,@("2022-02", 99.3)
    ,@("2022-03", 98.5)
    ,@("2022-04", 98.90000000000001)
    ,@("2022-05", 99.5)
    ,@("2022-06", 99.90000000000001)
    ,@("2022-07", 95.3)
    ,@("2022-08", 95.7)
    ,@("2022-09", 100.3)
    ,@("2023-01", 100.6)
    ,@("2023-02", 100.6)
    ,@("2023-03", 100.8)
    ,@("2023-04", 101.1)
    ,@("2023-05", 100.9)
    ,@("2023-06", 100.9)
    ,@("2023-07", 100.8)
)

# Original sheet only had 49 rows (2 through 49); the edit re-orders the
# existing months and appends new months through 2023-07, growing the used
# range to A1:D68.
$originalLastRow = 49
$newLastRow = 1 + $data.Count

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 4).Value = $item[1]
    $row++
}

# Rows beyond the original A1:D49 extent need blank B/C cells created (the
# sheet keeps explicit - empty - placeholders in columns B and C alongside
# the A/D data). Copying an existing blank B/C pair stamps the same blank
# cells onto the newly-used rows without touching their formatting/style.
if ($newLastRow -gt $originalLastRow) {
    $ws.Range("B2:C2").Copy($ws.Range("B" + ($originalLastRow + 1) + ":C" + $newLastRow))
}
